$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "code"
$ws.Range("B1").Value = "name"
$ws.Range("C1").Value = "lang_code"
$ws.Range("D1").Value = "is_active"

# Row 2 - Masculine
$ws.Range("A2").Value = "MLE"
$ws.Range("B2").Value = "MASCULIN"
$ws.Range("C2").Value = "fra"
$ws.Range("D2").Value = $true

# Row 3 - Feminine
$ws.Range("A3").Value = "FLE"
$ws.Range("B3").Value = "FEMININ"
$ws.Range("C3").Value = "fra"
$ws.Range("D3").Value = $true

# Update selection to match saved cursor position
$ws.Range("C9").Select()
